$wb = $excel.ActiveWorkbook

# 1. Update "Last Updated" timestamp on the Metadata sheet
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "05 Nov 2025, 12:02 PM"

# 2. Insert a new top row of data on the "Stock List" sheet.
#    A new ticker (CAPTRU-RE1) is inserted at row 2, pushing all the
#    existing rows down by one; the sheet keeps the same number of rows,
#    so the last existing row (77 -> originally 76, TRAVELFOOD) falls off.
$wsStock = $wb.Worksheets.Item("Stock List")

# Row number of the last populated row before the insert (data starts at row 2)
$lastRow = $wsStock.Cells.Item($wsStock.Rows.Count, 2).End(-4162).Row

# Insert a blank row at row 2, shifting everything below it (including the
# last row) down by one
$wsStock.Rows.Item(2).Insert(-4121)

# The inserted row inherits formatting from the row above (the bold header);
# reset it back to the plain/default style used by the rest of the data rows
$wsStock.Rows.Item(2).ClearFormats()

# Populate the newly inserted row 2 with the new stock data
$wsStock.Cells.Item(2, 1).Value = "📋"
$wsStock.Cells.Item(2, 2).Value = "CAPTRU-RE1"
$wsStock.Cells.Item(2, 3).Value = "CAPTRU-RE1"
$wsStock.Cells.Item(2, 4).Value = 5.67
$wsStock.Cells.Item(2, 5).Value = -11.9565
$wsStock.Cells.Item(2, 6).Value = "N/A"
$wsStock.Cells.Item(2, 7).Value = "N/A"
$wsStock.Cells.Item(2, 8).Value = 0

# The old last row of data is now one row further down (lastRow + 1); delete
# it so the sheet keeps the same total number of rows as before the edit
$wsStock.Rows.Item($lastRow + 1).Delete(-4162)
